# Applies the scheduled-runner price/profit refresh to the Leve profit sheets.
# Each block updates the currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N)
# for one leve row, recomputed from refreshed market-board data.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 52
$ws_ALC.Range("H52").Value = 764
$ws_ALC.Range("I52").Value = 764
$ws_ALC.Range("K52").Value = 2292
$ws_ALC.Range("M52").Value = -2132
# ALC row 68
$ws_ALC.Range("H68").Value = 0
$ws_ALC.Range("J68").Value = 0
$ws_ALC.Range("L68").Value = 0
$ws_ALC.Range("N68").ClearContents()
# ALC row 71
$ws_ALC.Range("H71").Value = 0
$ws_ALC.Range("J71").Value = 0
$ws_ALC.Range("L71").Value = 0
$ws_ALC.Range("N71").ClearContents()
# ALC row 97
$ws_ALC.Range("H97").Value = 4837.375
$ws_ALC.Range("I97").Value = 5000
$ws_ALC.Range("J97").Value = 4814.143
$ws_ALC.Range("K97").Value = 15000
$ws_ALC.Range("L97").Value = 14442.429
$ws_ALC.Range("M97").Value = -14504
$ws_ALC.Range("N97").Value = -15434.429
# ALC row 99
$ws_ALC.Range("H99").Value = 5888.5
$ws_ALC.Range("I99").Value = 219.5
$ws_ALC.Range("K99").Value = 658.5
$ws_ALC.Range("M99").Value = 839.5
# ALC row 137
$ws_ALC.Range("H137").Value = 2705.0908
$ws_ALC.Range("I137").Value = 2096.625
$ws_ALC.Range("J137").Value = 4327.6665
$ws_ALC.Range("K137").Value = 6289.875
$ws_ALC.Range("L137").Value = 12982.9995
$ws_ALC.Range("M137").Value = -3739.875
$ws_ALC.Range("N137").Value = -18082.9995
$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 46
$ws_ARM.Range("H46").Value = 27496.875
$ws_ARM.Range("I46").Value = 30000
$ws_ARM.Range("J46").Value = 27139.285
$ws_ARM.Range("K46").Value = 30000
$ws_ARM.Range("L46").Value = 27139.285
$ws_ARM.Range("M46").Value = -29681
$ws_ARM.Range("N46").Value = -27777.285
# ARM row 132
$ws_ARM.Range("H132").Value = 2384360.5
$ws_ARM.Range("I132").Value = 3470
$ws_ARM.Range("K132").Value = 10410
$ws_ARM.Range("M132").Value = -7880
$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws_BSM.Range("H105").Value = 618331.6
$ws_BSM.Range("I105").Value = 760650.9
$ws_BSM.Range("K105").Value = 760650.9
$ws_BSM.Range("M105").Value = -758903.9
$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws_CRP.Range("H7").Value = 176.5
$ws_CRP.Range("I7").Value = 243.33333
$ws_CRP.Range("K7").Value = 243.33333
$ws_CRP.Range("M7").Value = -130.33333
# CRP row 31
$ws_CRP.Range("H31").Value = 45456870
$ws_CRP.Range("I31").Value = 52633530
$ws_CRP.Range("J31").Value = 4718
$ws_CRP.Range("K31").Value = 52633530
$ws_CRP.Range("L31").Value = 4718
$ws_CRP.Range("M31").Value = -52633235
$ws_CRP.Range("N31").Value = -5308
# CRP row 34
$ws_CRP.Range("H34").Value = 45456870
$ws_CRP.Range("I34").Value = 52633530
$ws_CRP.Range("J34").Value = 4718
$ws_CRP.Range("K34").Value = 52633530
$ws_CRP.Range("L34").Value = 4718
$ws_CRP.Range("M34").Value = -52633328
$ws_CRP.Range("N34").Value = -5122
$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws_CUL.Range("H11").Value = 200009390
$ws_CUL.Range("I11").Value = 500000800
$ws_CUL.Range("K11").Value = 1500002400
$ws_CUL.Range("M11").Value = -1500002260
# CUL row 21
$ws_CUL.Range("H21").Value = 7856.7144
$ws_CUL.Range("I21").Value = 166.25
$ws_CUL.Range("J21").Value = 18110.666
$ws_CUL.Range("K21").Value = 498.75
$ws_CUL.Range("L21").Value = 54331.99800000001
$ws_CUL.Range("M21").Value = -325.75
$ws_CUL.Range("N21").Value = -54677.99800000001
# CUL row 22
$ws_CUL.Range("H22").Value = 12577.667
$ws_CUL.Range("I22").Value = 2200
$ws_CUL.Range("K22").Value = 6600
$ws_CUL.Range("M22").Value = -6431
# CUL row 25
$ws_CUL.Range("H25").Value = 33333
$ws_CUL.Range("I25").Value = 0
$ws_CUL.Range("K25").Value = 0
$ws_CUL.Range("M25").ClearContents()
# CUL row 27
$ws_CUL.Range("H27").Value = 12577.667
$ws_CUL.Range("I27").Value = 2200
$ws_CUL.Range("K27").Value = 6600
$ws_CUL.Range("M27").Value = -6498
# CUL row 30
$ws_CUL.Range("H30").Value = 33333
$ws_CUL.Range("I30").Value = 0
$ws_CUL.Range("K30").Value = 0
$ws_CUL.Range("M30").ClearContents()
# CUL row 35
$ws_CUL.Range("H35").Value = 9198.25
$ws_CUL.Range("I35").Value = 1305
$ws_CUL.Range("J35").Value = 17091.5
$ws_CUL.Range("K35").Value = 3915
$ws_CUL.Range("L35").Value = 51274.5
$ws_CUL.Range("M35").Value = -3627
$ws_CUL.Range("N35").Value = -51850.5
# CUL row 39
$ws_CUL.Range("H39").Value = 12628.875
$ws_CUL.Range("J39").Value = 14147.286
$ws_CUL.Range("L39").Value = 42441.858
$ws_CUL.Range("N39").Value = -43029.858
# CUL row 40
$ws_CUL.Range("H40").Value = 220.41667
$ws_CUL.Range("J40").Value = 308.2857
$ws_CUL.Range("L40").Value = 1233.1428
$ws_CUL.Range("N40").Value = -1371.1428
# CUL row 42
$ws_CUL.Range("H42").Value = 333345250
$ws_CUL.Range("I42").Value = 500001200
$ws_CUL.Range("K42").Value = 1500003600
$ws_CUL.Range("M42").Value = -1500003066
# CUL row 43
$ws_CUL.Range("H43").Value = 27008.25
$ws_CUL.Range("J43").Value = 27008.25
$ws_CUL.Range("L43").Value = 81024.75
$ws_CUL.Range("N43").Value = -81252.75
# CUL row 46
$ws_CUL.Range("H46").Value = 11474
$ws_CUL.Range("I46").Value = 890
$ws_CUL.Range("J46").Value = 16766
$ws_CUL.Range("K46").Value = 2670
$ws_CUL.Range("L46").Value = 50298
$ws_CUL.Range("M46").Value = -2579
$ws_CUL.Range("N46").Value = -50480
# CUL row 49
$ws_CUL.Range("H49").Value = 11186
$ws_CUL.Range("I49").Value = 112.5
$ws_CUL.Range("J49").Value = 33333
$ws_CUL.Range("K49").Value = 337.5
$ws_CUL.Range("L49").Value = 99999
$ws_CUL.Range("M49").Value = -181.5
$ws_CUL.Range("N49").Value = -100311
# CUL row 50
$ws_CUL.Range("H50").Value = 486.8
$ws_CUL.Range("I50").Value = 486.8
$ws_CUL.Range("J50").Value = 0
$ws_CUL.Range("K50").Value = 1460.4
$ws_CUL.Range("L50").Value = 0
$ws_CUL.Range("M50").Value = -979.4000000000001
$ws_CUL.Range("N50").ClearContents()
# CUL row 53
$ws_CUL.Range("H53").Value = 486.8
$ws_CUL.Range("I53").Value = 486.8
$ws_CUL.Range("J53").Value = 0
$ws_CUL.Range("K53").Value = 1460.4
$ws_CUL.Range("L53").Value = 0
$ws_CUL.Range("M53").Value = -979.4000000000001
$ws_CUL.Range("N53").ClearContents()
# CUL row 57
$ws_CUL.Range("H57").Value = 14665.6
$ws_CUL.Range("I57").Value = 9665.333000000001
$ws_CUL.Range("J57").Value = 22166
$ws_CUL.Range("K57").Value = 28995.999
$ws_CUL.Range("L57").Value = 66498
$ws_CUL.Range("M57").Value = -28436.999
$ws_CUL.Range("N57").Value = -67616
# CUL row 58
$ws_CUL.Range("H58").Value = 10482.308
$ws_CUL.Range("I58").Value = 8358
$ws_CUL.Range("K58").Value = 25074
$ws_CUL.Range("M58").Value = -24946
# CUL row 59
$ws_CUL.Range("H59").Value = 15110.667
$ws_CUL.Range("J59").Value = 22166
$ws_CUL.Range("L59").Value = 66498
$ws_CUL.Range("N59").Value = -67578
# CUL row 60
$ws_CUL.Range("H60").Value = 4749.4287
$ws_CUL.Range("J60").Value = 12802.4
$ws_CUL.Range("L60").Value = 38407.2
$ws_CUL.Range("N60").Value = -38909.2
# CUL row 93
$ws_CUL.Range("H93").Value = 11388.167
$ws_CUL.Range("J93").Value = 16082.25
$ws_CUL.Range("L93").Value = 48246.75
$ws_CUL.Range("N93").Value = -51990.75
# CUL row 122
$ws_CUL.Range("H122").Value = 11875.759
$ws_CUL.Range("J122").Value = 457.2
$ws_CUL.Range("L122").Value = 4114.8
$ws_CUL.Range("N122").Value = -9014.799999999999
# CUL row 132
$ws_CUL.Range("H132").Value = 1831.7727
$ws_CUL.Range("I132").Value = 1400.5294
$ws_CUL.Range("J132").Value = 3298
$ws_CUL.Range("K132").Value = 12604.7646
$ws_CUL.Range("L132").Value = 29682
$ws_CUL.Range("M132").Value = -10074.7646
$ws_CUL.Range("N132").Value = -34742
# CUL row 134
$ws_CUL.Range("H134").Value = 7659.3184
$ws_CUL.Range("I134").Value = 2509.5557
$ws_CUL.Range("K134").Value = 7528.6671
$ws_CUL.Range("M134").Value = -2458.6671
# CUL row 139
$ws_CUL.Range("H139").Value = 4483.1333
$ws_CUL.Range("I139").Value = 1814.3334
$ws_CUL.Range("K139").Value = 5443.0002
$ws_CUL.Range("M139").Value = -303.0002000000004
$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 7506.543
$ws_LTW.Range("I7").Value = 7270.864
$ws_LTW.Range("J7").Value = 7905.385
$ws_LTW.Range("K7").Value = 7270.864
$ws_LTW.Range("L7").Value = 7905.385
$ws_LTW.Range("M7").Value = -7158.864
$ws_LTW.Range("N7").Value = -8129.385
# LTW row 16
$ws_LTW.Range("H16").Value = 2779.8667
$ws_LTW.Range("J16").Value = 2498.5
$ws_LTW.Range("L16").Value = 2498.5
$ws_LTW.Range("N16").Value = -2838.5
# LTW row 68
$ws_LTW.Range("H68").Value = 2606892
$ws_LTW.Range("I68").Value = 5209915
$ws_LTW.Range("J68").Value = 3868.875
$ws_LTW.Range("K68").Value = 5209915
$ws_LTW.Range("L68").Value = 3868.875
$ws_LTW.Range("M68").Value = -5209166
$ws_LTW.Range("N68").Value = -5366.875
# LTW row 71
$ws_LTW.Range("H71").Value = 2606892
$ws_LTW.Range("I71").Value = 5209915
$ws_LTW.Range("J71").Value = 3868.875
$ws_LTW.Range("K71").Value = 26049575
$ws_LTW.Range("L71").Value = 19344.375
$ws_LTW.Range("M71").Value = -26045831
$ws_LTW.Range("N71").Value = -26832.375
# LTW row 119
$ws_LTW.Range("H119").Value = 85000
$ws_LTW.Range("J119").Value = 85000
$ws_LTW.Range("L119").Value = 85000
$ws_LTW.Range("N119").Value = -94676
# LTW row 126
$ws_LTW.Range("H126").Value = 7506.543
$ws_LTW.Range("I126").Value = 7270.864
$ws_LTW.Range("J126").Value = 7905.385
$ws_LTW.Range("K126").Value = 21812.592
$ws_LTW.Range("L126").Value = 23716.155
$ws_LTW.Range("M126").Value = -19342.592
$ws_LTW.Range("N126").Value = -28656.155
$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 126
$ws_WVR.Range("H126").Value = 6480.5
$ws_WVR.Range("I126").Value = 5283.769
$ws_WVR.Range("K126").Value = 15851.307
$ws_WVR.Range("M126").Value = -13381.307
# WVR row 136
$ws_WVR.Range("H136").Value = 220416.77
$ws_WVR.Range("I136").Value = 3392.3157
$ws_WVR.Range("J136").Value = 1251282.9
$ws_WVR.Range("K136").Value = 10176.9471
$ws_WVR.Range("L136").Value = 3753848.7
$ws_WVR.Range("M136").Value = -7626.947100000001
$ws_WVR.Range("N136").Value = -3758948.7
